# Refresh the crypto price/volume snapshot (Price column D, Volume(1h) column E)
# to match the latest scrape. Values are written with a leading apostrophe so
# they stay literal text (matching the sheet's existing inline-string cells)
# instead of being auto-converted to numbers/percentages by Excel.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'301.13"
$ws.Range("E2").Value = "'0.51%"
$ws.Range("D3").Value = "'32.15"
$ws.Range("E3").Value = "'1.43%"
$ws.Range("D4").Value = "'4.985"
$ws.Range("E4").Value = "'-3.08%"
$ws.Range("D5").Value = "'0.07892"
$ws.Range("E5").Value = "'-2.70%"
$ws.Range("D6").Value = "'2.097"
$ws.Range("E6").Value = "'-19.15%"
$ws.Range("D7").Value = "'7.796"
$ws.Range("E7").Value = "'0.12%"
$ws.Range("D8").Value = "'3.829"
$ws.Range("E8").Value = "'-2.10%"
$ws.Range("D9").Value = "'0.9283"
$ws.Range("E9").Value = "'-0.59%"
$ws.Range("D10").Value = "'0.1751"
$ws.Range("E10").Value = "'-0.26%"
$ws.Range("D11").Value = "'0.07957"
$ws.Range("E11").Value = "'7.99%"
$ws.Range("D12").Value = "'0.08673"
$ws.Range("E12").Value = "'-2.33%"
$ws.Range("D13").Value = "'0.03107"
$ws.Range("E13").Value = "'2.60%"
$ws.Range("D14").Value = "'0.1001"
$ws.Range("E14").Value = "'0.05%"
$ws.Range("D15").Value = "'0.001518"
$ws.Range("E15").Value = "'-0.28%"
$ws.Range("D16").Value = "'0.005976"
$ws.Range("E16").Value = "'3.36%"
$ws.Range("E17").Value = "'2,100.26%"
$ws.Range("E18").Value = "'-2.80%"
$ws.Range("E19").Value = "'-0.50%"
$ws.Range("D20").Value = "'0.3287"
$ws.Range("E21").Value = "'-2.21%"
$ws.Range("E22").Value = "'2.82%"
$ws.Range("D24").Value = "'0.04609"
$ws.Range("E24").Value = "'-0.24%"
$ws.Range("D25").Value = "'0.001236"
$ws.Range("E25").Value = "'-0.26%"
$ws.Range("D26").Value = "'0.004447"
$ws.Range("E26").Value = "'-1.65%"
$ws.Range("E27").Value = "'4.23%"
$ws.Range("E39").Value = "'-2.68%"
$ws.Range("D40").Value = "'0.04764"
$ws.Range("E40").Value = "'3.52%"
$ws.Range("D41").Value = "'0.007441"
$ws.Range("E41").Value = "'7.94%"
$ws.Range("D42").Value = "'0.1359"
$ws.Range("E42").Value = "'-1.15%"
$ws.Range("D43").Value = "'0.002274"
$ws.Range("E43").Value = "'6.01%"
$ws.Range("D44").Value = "'0.01131"
$ws.Range("E44").Value = "'9.01%"
$ws.Range("D45").Value = "'0.00005984"
$ws.Range("E45").Value = "'-3.50%"
$ws.Range("D46").Value = "'0.00000000750"
$ws.Range("E46").Value = "'-0.02%"
$ws.Range("D47").Value = "'0.003392"
$ws.Range("E47").Value = "'-59.61%"
$ws.Range("D48").Value = "'0.8234"
$ws.Range("E48").Value = "'10.02%"
$ws.Range("D49").Value = "'0.00002101"
$ws.Range("E49").Value = "'-0.02%"
$ws.Range("D50").Value = "'0.0002001"
$ws.Range("E50").Value = "'-0.02%"
